$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price column (D) for affected rows so numeric-looking
# strings (e.g. "0.999", "528.63") are preserved as text, matching the source data.
$ws.Range("D2:D10").NumberFormat = "@"
$ws.Range("D12:D21").NumberFormat = "@"
$ws.Range("D23:D25").NumberFormat = "@"
$ws.Range("D27:D34").NumberFormat = "@"
$ws.Range("D36:D51").NumberFormat = "@"

# Apply updated values row by row
$ws.Range("D2").Value = "58.397.11"
$ws.Range("E2").Value = "  -2.70%  "

$ws.Range("D3").Value = "2.453.82"
$ws.Range("E3").Value = "  -3.72%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Value = "528.63"
$ws.Range("E5").Value = "  -1.90%  "

$ws.Range("D6").Value = "133.94"
$ws.Range("E6").Value = "  -7.14%  "

$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("D8").Value = "0.555"
$ws.Range("E8").Value = "  -2.83%  "

$ws.Range("D9").Value = "2.460.49"
$ws.Range("E9").Value = "  -4.13%  "

$ws.Range("D10").Value = "0.0986"
$ws.Range("E10").Value = "  -2.68%  "

$ws.Range("E11").Value = "  -0.26%  "

$ws.Range("D12").Value = "5.30"
$ws.Range("E12").Value = "  -3.09%  "

$ws.Range("D13").Value = "0.343"
$ws.Range("E13").Value = "  -5.43%  "

$ws.Range("D14").Value = "2.886.47"
$ws.Range("E14").Value = "  -3.65%  "

$ws.Range("D15").Value = "58.222.65"
$ws.Range("E15").Value = "  -2.91%  "

$ws.Range("D16").Value = "22.59"
$ws.Range("E16").Value = "  -6.20%  "

$ws.Range("D17").Value = "0.0000139"
$ws.Range("E17").Value = "  -3.56%  "

$ws.Range("D18").Value = "2.455.51"
$ws.Range("E18").Value = "  -4.29%  "

$ws.Range("D19").Value = "10.72"
$ws.Range("E19").Value = "  -4.79%  "

$ws.Range("D20").Value = "4.19"
$ws.Range("E20").Value = "  -3.42%  "

$ws.Range("D21").Value = "321.11"
$ws.Range("E21").Value = "  -1.83%  "

$ws.Range("E22").Value = "  -0.25%  "

$ws.Range("D23").Value = "5.73"
$ws.Range("E23").Value = "  -3.59%  "

$ws.Range("D24").Value = "62.52"
$ws.Range("E24").Value = "  -0.95%  "

$ws.Range("D25").Value = "0.407"
$ws.Range("E25").Value = "  -6.31%  "

$ws.Range("E26").Value = "  -1.56%  "

$ws.Range("D27").Value = "0.982"
$ws.Range("E27").Value = "  -1.34%  "

$ws.Range("D28").Value = "7.45"
$ws.Range("E28").Value = "  -7.11%  "

$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0749"
$ws.Range("E29").Value = "  -5.90%  "

$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").Value = "6.48"
$ws.Range("E30").Value = "  -8.03%  "

$ws.Range("D31").Value = "1.75"
$ws.Range("E31").Value = "  -3.50%  "

$ws.Range("D32").Value = "163.49"
$ws.Range("E32").Value = "  -1.07%  "

$ws.Range("D33").Value = "0.998"
$ws.Range("E33").Value = "  +0.01%  "

$ws.Range("D34").Value = "1.08"
$ws.Range("E34").Value = "  -8.44%  "

$ws.Range("E35").Value = "  -8.07%  "

$ws.Range("D36").Value = "18.21"
$ws.Range("E36").Value = "  -2.80%  "

$ws.Range("D37").Value = "4.03"
$ws.Range("E37").Value = "  -8.55%  "

$ws.Range("D38").Value = "1.54"
$ws.Range("E38").Value = "  -5.62%  "

$ws.Range("D39").Value = "36.34"
$ws.Range("E39").Value = "  -1.77%  "

$ws.Range("B40").Value = "SuiNetwork"
$ws.Range("C40").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D40").Value = "0.795"
$ws.Range("E40").Value = "  -5.09%  "

$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "3.54"
$ws.Range("E41").Value = "  -5.12%  "

$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "275.49"
$ws.Range("E42").Value = "  -8.31%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "0.996"
$ws.Range("E43").Value = "  +0.24%  "

$ws.Range("D44").Value = "5.03"
$ws.Range("E44").Value = "  -9.90%  "

$ws.Range("D45").Value = "10.84"
$ws.Range("E45").Value = "  +0.15%  "

$ws.Range("D46").Value = "0.585"
$ws.Range("E46").Value = "  -3.80%  "

$ws.Range("D47").Value = "0.0921"
$ws.Range("E47").Value = "  -1.70%  "

$ws.Range("D48").Value = "120.48"
$ws.Range("E48").Value = "  -5.14%  "

$ws.Range("D49").Value = "0.0506"
$ws.Range("E49").Value = "  -2.58%  "

$ws.Range("D50").Value = "0.0217"
$ws.Range("E50").Value = "  -5.13%  "

$ws.Range("D51").Value = "17.05"
$ws.Range("E51").Value = "  -6.42%  "

